$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'24.699.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.71%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'1.678.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.61%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.42%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'313.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.19%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("E6").Value = "'  +0.33%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.3928"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.64%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.3973"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.09%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'1.004"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.41%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'51.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.74%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'1.410"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -5.00%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.08637"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.01%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'25.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -4.57%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'7.334"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.12%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'7.801"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.01%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'0.00001320"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.90%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'1.677.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.89%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'93.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.29%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'0.07087"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.95%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'20.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.94%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'7.098"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.64%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("E22").Value = "'  +0.36%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'13.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.73%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'24.704.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.67%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'2.367"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.40%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'23.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.79%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'2.767"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -7.39%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'163.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.03%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("B29").Value = "'HuobiToken"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'5.766"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -6.63%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("B30").Value = "'BitcoinCash"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'149.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.29%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'7.855"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -7.52%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'2.403"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +6.18%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'1.862.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.29%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'0.08444"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.93%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'0.03077"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.00%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'6.958"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.11%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'1.005"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.96%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'0.2794"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.13%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.09490"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.48%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'10.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.43%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.7933"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -5.73%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'1.486"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.96%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'13.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.51%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'16.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -5.05%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'0.7138"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -3.83%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'2.566"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.55%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'4.176"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.56%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.08687"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +4.02%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "'  +0.31%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'1.340"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.58%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'137.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.52%  "
$ws.Range("E51").Style = "Normal"
